$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Qmin (B) and Qmax (C) for rows 58-68 in the E772 table
for ($row = 58; $row -le 68; $row++) {
    $ws.Cells.Item($row, 2).Value = 11
    $ws.Cells.Item($row, 3).Value = 12
}
